$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "307.53"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-6.24%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "40.27"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-9.22%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.055"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-5.23%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07776"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-7.10%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.324"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-1.86%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.641"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-14.96%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9087"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-6.47%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1019"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-10.07%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1747"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-8.28%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09001"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-6.72%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04431"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-3.71%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.087"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-16.75%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1057"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.34%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001278"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.88%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005978"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "3.67%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.355"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.95%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.59%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3367"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.24%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1386"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.21%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2667"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "3.54%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04170"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.43%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001214"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.48%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004081"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-7.41%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001229"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-5.37%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003000"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.75%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02405"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-11.37%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05183"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-7.43%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007990"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.76%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1320"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-6.66%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007495"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "2.62%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001989"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-5.99%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008081"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-6.98%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3342"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-4.89%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006747"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-2.22%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000756"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.81%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003321"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-4.85%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004129"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "17.02%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002116"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.81%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002015"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.81%"
